$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet has weekly rows appended for the "Camote" variedad at Macroferia
# Regional de Talca. Two new records (1a (guarda) and 2a (guarda), both dated
# 45132) are inserted right after the existing row for 2024 and before the
# data that already starts at row 418, pushing everything from the old row
# 418 down by two rows (old 418 -> new 420, ..., old 460 -> new 462).
$ws.Range("A418:R419").Insert()

# New row 418: 1a (guarda)
$ws.Range("A418").Value = 5
$ws.Range("B418").Value = "Macroferia Regional de Talca"
$ws.Range("C418").Value = "Maule"
$ws.Range("D418").Value = 45132
$ws.Range("E418").Value = 7
$ws.Range("F418").Value = 100112045
$ws.Range("G418").Value = "Zapallo"
$ws.Range("H418").Value = "Camote"
$ws.Range("I418").Value = "1a (guarda)"
$ws.Range("J418").Value = 500
$ws.Range("K418").Value = 350
$ws.Range("L418").Value = 350
$ws.Range("M418").Value = 350
$ws.Range("N418").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O418").Value = "Región del Maule"
$ws.Range("P418").Value = 350
$ws.Range("Q418").Value = 1
$ws.Range("R418").Value = "Hortaliza"

# New row 419: 2a (guarda)
$ws.Range("A419").Value = 5
$ws.Range("B419").Value = "Macroferia Regional de Talca"
$ws.Range("C419").Value = "Maule"
$ws.Range("D419").Value = 45132
$ws.Range("E419").Value = 7
$ws.Range("F419").Value = 100112045
$ws.Range("G419").Value = "Zapallo"
$ws.Range("H419").Value = "Camote"
$ws.Range("I419").Value = "2a (guarda)"
$ws.Range("J419").Value = 500
$ws.Range("K419").Value = 280
$ws.Range("L419").Value = 280
$ws.Range("M419").Value = 280
$ws.Range("N419").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O419").Value = "Región del Maule"
$ws.Range("P419").Value = 280
$ws.Range("Q419").Value = 1
$ws.Range("R419").Value = "Hortaliza"
